$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 27714.428
$ws.Range("I63").Value = 8200
$ws.Range("J63").Value = 29215.54
$ws.Range("K63").Value = 8200
$ws.Range("L63").Value = 29215.54
$ws.Range("M63").Value = -7576
$ws.Range("N63").Value = -30463.54

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 27714.428
$ws.Range("I66").Value = 8200
$ws.Range("J66").Value = 29215.54
$ws.Range("K66").Value = 24600
$ws.Range("L66").Value = 87646.62
$ws.Range("M66").Value = -21480
$ws.Range("N66").Value = -93886.62

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3773.1667
$ws.Range("I116").Value = 4305.579
$ws.Range("J116").Value = 1750
$ws.Range("K116").Value = 4305.579
$ws.Range("L116").Value = 1750
$ws.Range("M116").Value = -863.5789999999997
$ws.Range("N116").Value = -8634

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1050.7407
$ws.Range("I121").Value = 614.1667
$ws.Range("J121").Value = 1175.4762
$ws.Range("K121").Value = 1842.5001
$ws.Range("L121").Value = 3526.4286
$ws.Range("M121").Value = -95.50009999999997
$ws.Range("N121").Value = -7020.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2752.1282
$ws.Range("I138").Value = 1657.0385
$ws.Range("J138").Value = 4942.3076
$ws.Range("K138").Value = 4971.1155
$ws.Range("L138").Value = 14826.9228
$ws.Range("M138").Value = 168.8845000000001
$ws.Range("N138").Value = -25106.9228

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3839.5952
$ws.Range("I141").Value = 1372.579
$ws.Range("J141").Value = 27276.25
$ws.Range("K141").Value = 4117.737
$ws.Range("L141").Value = 81828.75
$ws.Range("M141").Value = 1062.263
$ws.Range("N141").Value = -92188.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 655.85187
$ws.Range("I2").Value = 542.36365
$ws.Range("J2").Value = 1155.2
$ws.Range("K2").Value = 542.36365
$ws.Range("L2").Value = 1155.2
$ws.Range("M2").Value = -429.36365
$ws.Range("N2").Value = -1381.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 41000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 41000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 41000
$ws.Range("N62").Value = -42248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 41000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 41000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 123000
$ws.Range("N65").Value = -129240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1556.3448
$ws.Range("I97").Value = 540.5
$ws.Range("J97").Value = 30000
$ws.Range("K97").Value = 540.5
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = -44.5
$ws.Range("N97").Value = -30992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N102").ClearContents()
$ws.Range("H102").Value = 1730
$ws.Range("I102").Value = 1730
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1730
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1852.25
$ws.Range("I110").Value = 1727.7646
$ws.Range("J110").Value = 2275.5
$ws.Range("K110").Value = 1727.7646
$ws.Range("L110").Value = 2275.5
$ws.Range("M110").Value = 317.2354
$ws.Range("N110").Value = -6365.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 655.85187
$ws.Range("I116").Value = 542.36365
$ws.Range("J116").Value = 1155.2
$ws.Range("K116").Value = 542.36365
$ws.Range("L116").Value = 1155.2
$ws.Range("M116").Value = 1751.63635
$ws.Range("N116").Value = -5743.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4520.711
$ws.Range("I132").Value = 5144.9033
$ws.Range("J132").Value = 3138.5715
$ws.Range("K132").Value = 15434.7099
$ws.Range("L132").Value = 9415.7145
$ws.Range("M132").Value = -12904.7099
$ws.Range("N132").Value = -14475.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 655.85187
$ws.Range("I3").Value = 542.36365
$ws.Range("J3").Value = 1155.2
$ws.Range("K3").Value = 542.36365
$ws.Range("L3").Value = 1155.2
$ws.Range("M3").Value = -428.36365
$ws.Range("N3").Value = -1383.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 304.125
$ws.Range("I64").Value = 151.75
$ws.Range("J64").Value = 354.91666
$ws.Range("K64").Value = 151.75
$ws.Range("L64").Value = 354.91666
$ws.Range("M64").Value = 73.25
$ws.Range("N64").Value = -804.91666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 304.125
$ws.Range("I67").Value = 151.75
$ws.Range("J67").Value = 354.91666
$ws.Range("K67").Value = 151.75
$ws.Range("L67").Value = 354.91666
$ws.Range("M67").Value = 628.25
$ws.Range("N67").Value = -1914.91666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1023.44446
$ws.Range("I99").Value = 888.8889
$ws.Range("J99").Value = 1158
$ws.Range("K99").Value = 888.8889
$ws.Range("L99").Value = 1158
$ws.Range("M99").Value = 609.1111
$ws.Range("N99").Value = -4154

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2399.282
$ws.Range("I105").Value = 1797.5
$ws.Range("J105").Value = 3032.7368
$ws.Range("K105").Value = 1797.5
$ws.Range("L105").Value = 3032.7368
$ws.Range("M105").Value = -50.5
$ws.Range("N105").Value = -6526.736800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1329.54
$ws.Range("I31").Value = 881.93335
$ws.Range("J31").Value = 1521.3715
$ws.Range("K31").Value = 881.93335
$ws.Range("L31").Value = 1521.3715
$ws.Range("M31").Value = -586.93335
$ws.Range("N31").Value = -2111.3715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1329.54
$ws.Range("I34").Value = 881.93335
$ws.Range("J34").Value = 1521.3715
$ws.Range("K34").Value = 881.93335
$ws.Range("L34").Value = 1521.3715
$ws.Range("M34").Value = -679.93335
$ws.Range("N34").Value = -1925.3715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 38469748
$ws.Range("I132").Value = 83344330
$ws.Range("J132").Value = 5824.2856
$ws.Range("K132").Value = 250032990
$ws.Range("L132").Value = 17472.8568
$ws.Range("M132").Value = -250030460
$ws.Range("N132").Value = -22532.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2481.9412
$ws.Range("I134").Value = 2555.4167
$ws.Range("J134").Value = 2305.6
$ws.Range("K134").Value = 7666.250100000001
$ws.Range("L134").Value = 6916.799999999999
$ws.Range("M134").Value = -5131.250100000001
$ws.Range("N134").Value = -11986.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 721.01495
$ws.Range("I107").Value = 330.47726
$ws.Range("J107").Value = 1468.1305
$ws.Range("K107").Value = 991.43178
$ws.Range("L107").Value = 4404.3915
$ws.Range("M107").Value = 928.56822
$ws.Range("N107").Value = -8244.3915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1108.3334
$ws.Range("I109").Value = 1575
$ws.Range("J109").Value = 875
$ws.Range("K109").Value = 4725
$ws.Range("L109").Value = 2625
$ws.Range("M109").Value = -3685
$ws.Range("N109").Value = -4705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -188

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4599.9287
$ws.Range("I70").Value = 4706.25
$ws.Range("J70").Value = 4458.1665
$ws.Range("K70").Value = 4706.25
$ws.Range("L70").Value = 4458.1665
$ws.Range("M70").Value = -4436.25
$ws.Range("N70").Value = -4998.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4599.9287
$ws.Range("I73").Value = 4706.25
$ws.Range("J73").Value = 4458.1665
$ws.Range("K73").Value = 4706.25
$ws.Range("L73").Value = 4458.1665
$ws.Range("M73").Value = -3770.25
$ws.Range("N73").Value = -6330.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4076.6316
$ws.Range("I132").Value = 4524.2896
$ws.Range("J132").Value = 3181.3157
$ws.Range("K132").Value = 13572.8688
$ws.Range("L132").Value = 9543.947100000001
$ws.Range("M132").Value = -11042.8688
$ws.Range("N132").Value = -14603.9471

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6001.6577
$ws.Range("I132").Value = 6821.6665
$ws.Range("J132").Value = 3988.9092
$ws.Range("K132").Value = 20464.9995
$ws.Range("L132").Value = 11966.7276
$ws.Range("M132").Value = -17934.9995
$ws.Range("N132").Value = -17026.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4647.4546
$ws.Range("I136").Value = 6614.778
$ws.Range("J136").Value = 2286.6667
$ws.Range("K136").Value = 19844.334
$ws.Range("L136").Value = 6860.000100000001
$ws.Range("M136").Value = -17294.334
$ws.Range("N136").Value = -11960.0001
